# Auto-generated edit script for cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.279.91"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.288.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.283.45"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.100"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("E12").Value = "  +1.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.63"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.696.14"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.201.21"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.298.65"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.19"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.39%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  -1.43%  "

$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.65"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  -2.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0725"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.08"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.85"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("E38").Value = "  -0.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.92"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.50"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "289.34"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0951"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.17%  "

$ws.Range("E48").Value = "  -1.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.94"
$ws.Range("D49").ClearFormats()

$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("E51").Value = "  +1.50%  "
